# Applies the "Added the code for common function" change to
# AddProductPage.data.xlsx: five new label/selector rows appended to the
# (only) worksheet, and the selection moved to the new last cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 144-148 (column A = element/function name, column B = value).
$ws.Range("A144").Value = "columnBusiness"
$ws.Range("B144").Value = "//body/div//section//div//div[3]//table//tr/td[1]"

$ws.Range("A145").Value = "columnName3"
$ws.Range("B145").Value = "//table[@id='list-Business']//div[contains(text(),'Name')]"

$ws.Range("A146").Value = "columnName1"
$ws.Range("B146").Value = "//table[@id='"

$ws.Range("A147").Value = "columnName2"
$ws.Range("B147").Value = "]//div[contains(text(),'Name')]"

$ws.Range("A148").Value = "columnName"
$ws.Range("B148").Value = "list-Business'"

# Match the saved selection/viewport from the commit (cursor left on the
# last filled cell of the newly-added block).
$ws.Range("C148").Select() | Out-Null
